$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# Split the run "not final and I will have to perform tests and check if
# this approach succeeds in generating successful creatures." into three
# runs, wrapping "final" with a grammar proofErr pair (gramStart/gramEnd).
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$target1 = "All of this is not final and I will have to perform tests and check if this approach succeeds in generating successful creatures."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $target1) {
        $body = '<w:body><w:p>' +
            '<w:r><w:t xml:space="preserve">All of this is </w:t></w:r>' +
            '<w:r><w:t xml:space="preserve">not </w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>final</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '<w:r><w:t xml:space="preserve"> and I will have to perform tests and check if this approach succeeds in generating successful creatures.</w:t></w:r>' +
            '</w:p></w:body>'
        $null = $p.Range.InsertXML($pkgHeader + $body + $pkgFooter)
        break
    }
}

# --- Change 2 -------------------------------------------------------------
# Insert a brand-new paragraph right after the "DAY 03 ... START UP UNITY!"
# paragraph (and before the trailing empty paragraph) describing the extra
# pseudocode work, with spelling (pseudocde) and grammar (bs) proofErr
# markers.
$target2 = "DAY 03: Pseudocode is good. Review it and see if anymore progress is needed. If not, START UP UNITY!"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $target2) {
        $body = '<w:body>' +
            '<w:p>' +
                '<w:r><w:t xml:space="preserve">DAY 03: Pseudocode is good. Review it and see if anymore </w:t></w:r>' +
                '<w:r><w:t>progress</w:t></w:r>' +
                '<w:r><w:t xml:space="preserve"> is needed. If not, START UP UNITY!</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
                '<w:r><w:t xml:space="preserve">Started working a bit more on </w:t></w:r>' +
                '<w:proofErr w:type="spellStart"/>' +
                '<w:r><w:t>pseudocde</w:t></w:r>' +
                '<w:proofErr w:type="spellEnd"/>' +
                '<w:r><w:t xml:space="preserve">. Mind is unclear so I will stop here, I believe most of my code is straight </w:t></w:r>' +
                '<w:proofErr w:type="gramStart"/>' +
                '<w:r><w:t>bs</w:t></w:r>' +
                '<w:proofErr w:type="gramEnd"/>' +
                '<w:r><w:t xml:space="preserve"> so I’ll do it with a clear head.</w:t></w:r>' +
            '</w:p>' +
            '</w:body>'
        $null = $p.Range.InsertXML($pkgHeader + $body + $pkgFooter)
        break
    }
}
